$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.358.80"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.245.92"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "247.67"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "75.96"
$ws.Range("E7").Value = "  +6.13%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "40.03"
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "7.22"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "2.587.48"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "14.86"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "2.264.71"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "42.362.77"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "0.0₃0979"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "71.42"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "231.08"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "2.15"
$ws.Range("E23").Value = "  -4.87%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("E25").Value = "  -4.81%  "
$ws.Range("D26").Value = "11.14"
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("D27").Value = "2.31"
$ws.Range("E27").Value = "  -5.03%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "6.90"
$ws.Range("E29").Value = "  +7.95%  "
$ws.Range("D30").Value = "168.08"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").Value = "20.52"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "0.0852"
$ws.Range("E32").Value = "  +7.18%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "31.07"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.120"
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").Value = "4.49"
$ws.Range("E36").Value = "  -5.43%  "
$ws.Range("D37").Value = "4.72"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "0.0297"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "12.96"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").Value = "5.91"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "116.64"
$ws.Range("E42").Value = "  +22.22%  "
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").Value = "60.02"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").Value = "8.75"
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "4.11"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "2.69"
$ws.Range("E51").Value = "  -1.78%  "
